$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed cryptocurrency Price (D) and Volume(1h) (E) figures for rows 2-51.
# D = $null means the displayed Price text is unchanged for that row.
$updates = @(
    @{ Row = 2; D = '42.045.75'; E = '  -0.87%  ' }
    @{ Row = 3; D = '2.237.97'; E = '  +0.07%  ' }
    @{ Row = 4; D = $null; E = '  -0.06%  ' }
    @{ Row = 5; D = '242.85'; E = '  -0.89%  ' }
    @{ Row = 6; D = $null; E = '  -0.72%  ' }
    @{ Row = 7; D = '74.19'; E = '  -0.12%  ' }
    @{ Row = 8; D = $null; E = '  +0.10%  ' }
    @{ Row = 9; D = $null; E = '  -3.74%  ' }
    @{ Row = 10; D = '42.14'; E = '  -2.21%  ' }
    @{ Row = 11; D = '0.0945'; E = '  -1.51%  ' }
    @{ Row = 12; D = $null; E = '  +0.17%  ' }
    @{ Row = 13; D = '6.91'; E = '  -2.92%  ' }
    @{ Row = 14; D = '2.572.35'; E = '  +0.18%  ' }
    @{ Row = 15; D = '14.35'; E = '  -0.79%  ' }
    @{ Row = 16; D = '0.837'; E = '  -1.78%  ' }
    @{ Row = 17; D = '2.251.76'; E = '  +1.30%  ' }
    @{ Row = 18; D = '41.952.72'; E = '  -0.64%  ' }
    @{ Row = 19; D = $null; E = '  -3.75%  ' }
    @{ Row = 20; D = '6.20'; E = '  +0.60%  ' }
    @{ Row = 21; D = '72.58'; E = '  +0.69%  ' }
    @{ Row = 22; D = '11.21'; E = '  +9.82%  ' }
    @{ Row = 23; D = '229.79'; E = '  -0.71%  ' }
    @{ Row = 24; D = $null; E = '  -6.35%  ' }
    @{ Row = 25; D = $null; E = '  +0.19%  ' }
    @{ Row = 26; D = '11.36'; E = '  -3.19%  ' }
    @{ Row = 27; D = $null; E = '  -0.59%  ' }
    @{ Row = 28; D = '2.28'; E = '  -0.85%  ' }
    @{ Row = 29; D = $null; E = '  -0.90%  ' }
    @{ Row = 30; D = '167.62'; E = '  +0.56%  ' }
    @{ Row = 31; D = '20.56'; E = '  -1.79%  ' }
    @{ Row = 32; D = '5.58'; E = '  -4.40%  ' }
    @{ Row = 33; D = '0.0799'; E = '  -0.92%  ' }
    @{ Row = 34; D = '30.02'; E = '  +1.16%  ' }
    @{ Row = 35; D = $null; E = '  -0.74%  ' }
    @{ Row = 36; D = $null; E = '  -6.45%  ' }
    @{ Row = 37; D = '4.27'; E = '  -4.47%  ' }
    @{ Row = 38; D = $null; E = '  -1.58%  ' }
    @{ Row = 39; D = '13.09'; E = '  -1.24%  ' }
    @{ Row = 40; D = $null; E = '  -2.05%  ' }
    @{ Row = 41; D = $null; E = '  +0.76%  ' }
    @{ Row = 42; D = '64.58'; E = '  +1.79%  ' }
    @{ Row = 43; D = '0.198'; E = '  -1.56%  ' }
    @{ Row = 44; D = '8.69'; E = '  -1.48%  ' }
    @{ Row = 45; D = '103.63'; E = '  -2.11%  ' }
    @{ Row = 46; D = $null; E = '  -1.98%  ' }
    @{ Row = 47; D = $null; E = '  -0.34%  ' }
    @{ Row = 48; D = $null; E = '  -0.90%  ' }
    @{ Row = 49; D = $null; E = '  -2.23%  ' }
    @{ Row = 50; D = $null; E = '  -2.03%  ' }
    @{ Row = 51; D = '2.447.75'; E = '  +0.05%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Range("D$($u.Row)")
        if ($u.D -match "^[+-]?[0-9]*\.?[0-9]+$") {
            # The string looks like a plain number (e.g. "242.85"); force Text
            # format first so Excel keeps it as a literal string instead of
            # silently converting it to a numeric value.
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    $ws.Range("E$($u.Row)").Value = $u.E
}
